# Adjusting period ranges for the four models
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for rows 2-21 (Year_Range, Keyword, Correlation)
$data = @(
    @("2008_2009", "inflation", -0.4892),
    @("2008_2009", "interest", -0.4244),
    @("2008_2009", "uncertain", -0.3641),
    @("2008_2009", "invest", 0.0982),
    @("2008_2009", "trade", -0.7977),
    @("2010_2019", "uncertain", -0.3662),
    @("2010_2019", "interest", 0.4306),
    @("2010_2019", "trade", 0.192),
    @("2010_2019", "invest", 0.4317),
    @("2010_2019", "inflation", 0.7954),
    @("2020_2021", "inflation", -0.3408),
    @("2020_2021", "interest", -0.2763),
    @("2020_2021", "invest", -0.2103),
    @("2020_2021", "trade", 0.0678),
    @("2020_2021", "uncertain", -0.5348000000000001),
    @("2022_2023", "inflation", -0.7396),
    @("2022_2023", "interest", -0.6764),
    @("2022_2023", "uncertain", -0.3411),
    @("2022_2023", "invest", 0.3824),
    @("2022_2023", "trade", -0.3993)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
